$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The engine's row-delete does not reliably shift single-cell merged ranges
# (e.g. "A42", a merge of just one cell) the way it shifts multi-cell merged
# ranges. Work around this by unmerging any single-cell merges below the row
# being removed, then re-merging them (at their new, shifted location) after
# the delete.
$ws.Range("A42").UnMerge()
$ws.Range("B42").UnMerge()

# Delete row 15 (duplicate "checksum" row), shifting all rows below up by one.
$ws.Rows.Item(15).Delete()

$ws.Range("A41").Merge()
$ws.Range("B41").Merge()
